$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 46070, -18.684, 915.3, "17.02.20261"),
    @(3, 46070.01041666666, -41.51, 1972.376, "17.02.20262"),
    @(4, 46070.02083333334, -50.586, 932.272, "17.02.20263"),
    @(5, 46070.03125, -49.329, 1013.184, "17.02.20264"),
    @(6, 46070.04166666666, -44.534, 907.5410000000001, "17.02.20265"),
    @(7, 46070.05208333334, -53.741, 1502.3, "17.02.20266"),
    @(8, 46070.0625, -38.394, 892.365, "17.02.20267"),
    @(9, 46070.07291666666, -34.768, 893.853, "17.02.20268"),
    @(10, 46070.08333333334, -56.729, 892.88, "17.02.20269"),
    @(11, 46070.09375, -76.905, 905.5410000000001, "17.02.202610"),
    @(12, 46070.10416666666, -100.529, 907.803, "17.02.202611"),
    @(13, 46070.11458333334, -89.99299999999999, 909.889, "17.02.202612"),
    @(14, 46070.125, -73.23099999999999, 999.636, "17.02.202613"),
    @(15, 46070.13541666666, -81.944, 913.837, "17.02.202614"),
    @(16, 46070.14583333334, -84.123, 908.192, "17.02.202615"),
    @(17, 46070.15625, -85.03100000000001, 904.309, "17.02.202616"),
    @(18, 46070.16666666666, -96.56, 1555.439, "17.02.202617"),
    @(19, 46070.17708333334, -141.228, 1105.107, "17.02.202618"),
    @(20, 46070.1875, -132.666, 909.206, "17.02.202619"),
    @(21, 46070.19791666666, -68.67100000000001, 914.396, "17.02.202620"),
    @(22, 46070.20833333334, -29.427, 935.986, "17.02.202621"),
    @(23, 46070.21875, -33.839, 934.99, "17.02.202622"),
    @(24, 46070.22916666666, -61.326, 1022.932, "17.02.202623"),
    @(25, 46070.23958333334, -60.959, 959.304, "17.02.202624"),
    @(26, 46070.25, -19.256, 979.927, "17.02.202625"),
    @(27, 46070.26041666666, -28.987, 783.98, "17.02.202626"),
    @(28, 46070.27083333334, -62.361, 679.864, "17.02.202627"),
    @(29, 46070.28125, -75.539, 632.8819999999999, "17.02.202628"),
    @(30, 46070.29166666666, -74.825, 3145.203, "17.02.202629"),
    @(31, 46070.30208333334, -43.735, 971.51, "17.02.202630"),
    @(32, 46070.3125, -29.543, 1085.708, "17.02.202631"),
    @(33, 46070.32291666666, 1.505, 371.791, "17.02.202632"),
    @(34, 46070.33333333334, 2.007, 311.912, "17.02.202633"),
    @(35, 46070.34375, 4.997, 400, "17.02.202634"),
    @(36, 46070.35416666666, -6.255, 772.449, "17.02.202635"),
    @(37, 46070.36458333334, 13.228, 358.542, "17.02.202636"),
    @(38, 46070.375, -46.39, 1139.139, "17.02.202637"),
    @(39, 46070.38541666666, -123.422, 4816.128, "17.02.202638"),
    @(40, 46070.39583333334, -88.05500000000001, 1074.685, "17.02.202639"),
    @(41, 46070.40625, -134.566, 1212.905, "17.02.202640"),
    @(42, 46070.41666666666, -84.65600000000001, 1022.454, "17.02.202641"),
    @(43, 46070.42708333334, -138.113, 1014.531, "17.02.202642"),
    @(44, 46070.4375, -157.433, 1017.91, "17.02.202643"),
    @(45, 46070.44791666666, -110.557, 1751.459, "17.02.202644"),
    @(46, 46070.45833333334, -57.106, 972.847, "17.02.202645"),
    @(47, 46070.46875, -35.038, 996, "17.02.202646"),
    @(48, 46070.47916666666, 21.732, 399.641, "17.02.202647"),
    @(49, 46070.48958333334, 0.571, 246.983, "17.02.202648"),
    @(50, 46070.5, 13.746, 255.807, "17.02.202649"),
    @(51, 46070.51041666666, 17.501, 226.85, "17.02.202650"),
    @(52, 46070.52083333334, -2.11, 746.621, "17.02.202651"),
    @(53, 46070.53125, 4.501, 332.579, "17.02.202652"),
    @(54, 46070.54166666666, -37.362, 813.612, "17.02.202653"),
    @(55, 46070.55208333334, -11.909, 759.895, "17.02.202654"),
    @(56, 46070.5625, -33.548, 804.039, "17.02.202655"),
    @(57, 46070.57291666666, -23.574, 862.037, "17.02.202656"),
    @(58, 46070.58333333334, -48.794, 3585.469, "17.02.202657"),
    @(59, 46070.59375, -33.952, 1332.565, "17.02.202658"),
    @(60, 46070.60416666666, 0.142, 711.474, "17.02.202659"),
    @(61, 46070.61458333334, -6.833, 857.192, "17.02.202660"),
    @(62, 46070.625, -16.244, 827.135, "17.02.202661"),
    @(63, 46070.63541666666, 10.374, 500, "17.02.202662"),
    @(64, 46070.64583333334, 13.811, 0, "17.02.202663"),
    @(65, 46070.65625, 47.666, 343.971, "17.02.202664"),
    @(66, 46070.66666666666, 23.511, 289.382, "17.02.202665"),
    @(67, 46070.67708333334, 33.435, 443.83, "17.02.202666"),
    @(68, 46070.6875, -4.306, 780.529, "17.02.202667"),
    @(69, 46070.69791666666, -7.068, 851.141, "17.02.202668"),
    @(70, 46070.70833333334, -18.782, 984.669, "17.02.202669"),
    @(71, 46070.71875, -2.32, 980.658, "17.02.202670"),
    @(72, 46070.72916666666, -8.085000000000001, 794.732, "17.02.202671"),
    @(73, 46070.73958333334, -1.536, 796.395, "17.02.202672"),
    @(74, 46070.75, -52.863, 807.912, "17.02.202673"),
    @(75, 46070.76041666666, -25.891, 882.883, "17.02.202674"),
    @(76, 46070.77083333334, -12.681, 769.962, "17.02.202675"),
    @(77, 46070.78125, 11.213, 277.807, "17.02.202676"),
    @(78, 46070.79166666666, -2.173, 599, "17.02.202677"),
    @(79, 46070.80208333334, -1.053, 709.898, "17.02.202678"),
    @(80, 46070.8125, 13.771, 258.394, "17.02.202679"),
    @(81, 46070.82291666666, -5.9, 867.422, "17.02.202680"),
    @(82, 46070.83333333334, -9.693, 832.479, "17.02.202681"),
    @(83, 46070.84375, -10.813, 790.203, "17.02.202682"),
    @(84, 46070.85416666666, 3.808, 305.855, "17.02.202683"),
    @(85, 46070.86458333334, -6.821, 668.811, "17.02.202684"),
    @(86, 46070.875, -54.196, 890.74, "17.02.202685"),
    @(87, 46070.88541666666, -49.146, 878.145, "17.02.202686"),
    @(88, 46070.89583333334, -30.171, 869.127, "17.02.202687"),
    @(89, 46070.90625, -32.237, 867.617, "17.02.202688"),
    @(90, 46070.91666666666, -56.711, 2295.746, "17.02.202689"),
    @(91, 46070.92708333334, -39.07, 1200.973, "17.02.202690"),
    @(92, 46070.9375, -35.945, 818.171, "17.02.202691"),
    @(93, 46070.94791666666, -30.614, 816.037, "17.02.202692"),
    @(94, 46070.95833333334, -90.93899999999999, 992.722, "17.02.202693"),
    @(95, 46070.96875, -64.286, 727.223, "17.02.202694"),
    @(96, 46070.97916666666, -63.779, 875.752, "17.02.202695"),
    @(97, 46070.98958333334, -61.613, 716.8339999999999, "17.02.202696"),
    @(98, 46071, -40.331, 507.022, "18.02.20261"),
    @(99, 46071.01041666666, -36.325, 510.951, "18.02.20262"),
    @(100, 46071.02083333334, -77.199, 548.38, "18.02.20263"),
    @(101, 46071.03125, -83.749, 507.036, "18.02.20264"),
    @(102, 46071.04166666666, -44.822, 542.023, "18.02.20265"),
    @(103, 46071.05208333334, -54.425, 539.59, "18.02.20266"),
    @(104, 46071.0625, -57.416, 588.087, "18.02.20267"),
    @(105, 46071.07291666666, -73.73099999999999, 507.676, "18.02.20268"),
    @(106, 46071.08333333334, -64.119, 533.218, "18.02.20269"),
    @(107, 46071.09375, -37.495, 511.908, "18.02.202610"),
    @(108, 46071.10416666666, -85.30200000000001, 510.078, "18.02.202611"),
    @(109, 46071.11458333334, -115.321, 510.03, "18.02.202612"),
    @(110, 46071.125, -68.864, 588.39, "18.02.202613"),
    @(111, 46071.13541666666, -76.62, 734.651, "18.02.202614"),
    @(112, 46071.14583333334, -113.079, 834.208, "18.02.202615"),
    @(113, 46071.15625, -143.817, 563.0119999999999, "18.02.202616"),
    @(114, 46071.16666666666, -68.045, 525.74, "18.02.202617"),
    @(115, 46071.17708333334, -89.97499999999999, 525.748, "18.02.202618"),
    @(116, 46071.1875, -135.882, 1105.406, "18.02.202619"),
    @(117, 46071.19791666666, -108.768, 627.355, "18.02.202620"),
    @(118, 46071.20833333334, -24.418, 830.521, "18.02.202621"),
    @(119, 46071.21875, -42.378, 834.886, "18.02.202622"),
    @(120, 46071.22916666666, -64.185, 826.25, "18.02.202623"),
    @(121, 46071.23958333334, -92.614, 821.918, "18.02.202624"),
    @(122, 46071.25, -9.272, 1000, "18.02.202625"),
    @(123, 46071.26041666666, -37.487, 999.599, "18.02.202626"),
    @(124, 46071.27083333334, -24.636, 998.8200000000001, "18.02.202627"),
    @(125, 46071.28125, -54.281, 999.968, "18.02.202628"),
    @(126, 46071.29166666666, -61.406, 0, "18.02.202629"),
    @(127, 46071.30208333334, 0, 0, "18.02.202630"),
    @(128, 46071.3125, 0, 0, "18.02.202631"),
    @(129, 46071.32291666666, 0, 0, "18.02.202632"),
    @(130, 46071.33333333334, 0, 0, "18.02.202633"),
    @(131, 46071.34375, 0, 0, "18.02.202634"),
    @(132, 46071.35416666666, 0, 0, "18.02.202635"),
    @(133, 46071.36458333334, 0, 0, "18.02.202636"),
    @(134, 46071.375, 0, 0, "18.02.202637"),
    @(135, 46071.38541666666, 0, 0, "18.02.202638"),
    @(136, 46071.39583333334, 0, 0, "18.02.202639"),
    @(137, 46071.40625, 0, 0, "18.02.202640"),
    @(138, 46071.41666666666, 0, 0, "18.02.202641"),
    @(139, 46071.42708333334, 0, 0, "18.02.202642"),
    @(140, 46071.4375, 0, 0, "18.02.202643"),
    @(141, 46071.44791666666, 0, 0, "18.02.202644"),
    @(142, 46071.45833333334, 0, 0, "18.02.202645"),
    @(143, 46071.46875, 0, 0, "18.02.202646"),
    @(144, 46071.47916666666, 0, 0, "18.02.202647"),
    @(145, 46071.48958333334, 0, 0, "18.02.202648"),
    @(146, 46071.5, 0, 0, "18.02.202649"),
    @(147, 46071.51041666666, 0, 0, "18.02.202650"),
    @(148, 46071.52083333334, 0, 0, "18.02.202651"),
    @(149, 46071.53125, 0, 0, "18.02.202652"),
    @(150, 46071.54166666666, 0, 0, "18.02.202653"),
    @(151, 46071.55208333334, 0, 0, "18.02.202654"),
    @(152, 46071.5625, 0, 0, "18.02.202655"),
    @(153, 46071.57291666666, 0, 0, "18.02.202656"),
    @(154, 46071.58333333334, 0, 0, "18.02.202657"),
    @(155, 46071.59375, 0, 0, "18.02.202658"),
    @(156, 46071.60416666666, 0, 0, "18.02.202659"),
    @(157, 46071.61458333334, 0, 0, "18.02.202660"),
    @(158, 46071.625, 0, 0, "18.02.202661"),
    @(159, 46071.63541666666, 0, 0, "18.02.202662"),
    @(160, 46071.64583333334, 0, 0, "18.02.202663"),
    @(161, 46071.65625, 0, 0, "18.02.202664"),
    @(162, 46071.66666666666, 0, 0, "18.02.202665"),
    @(163, 46071.67708333334, 0, 0, "18.02.202666"),
    @(164, 46071.6875, 0, 0, "18.02.202667"),
    @(165, 46071.69791666666, 0, 0, "18.02.202668"),
    @(166, 46071.70833333334, 0, 0, "18.02.202669"),
    @(167, 46071.71875, 0, 0, "18.02.202670"),
    @(168, 46071.72916666666, 0, 0, "18.02.202671"),
    @(169, 46071.73958333334, 0, 0, "18.02.202672"),
    @(170, 46071.75, 0, 0, "18.02.202673"),
    @(171, 46071.76041666666, 0, 0, "18.02.202674"),
    @(172, 46071.77083333334, 0, 0, "18.02.202675"),
    @(173, 46071.78125, 0, 0, "18.02.202676"),
    @(174, 46071.79166666666, 0, 0, "18.02.202677"),
    @(175, 46071.80208333334, 0, 0, "18.02.202678"),
    @(176, 46071.8125, 0, 0, "18.02.202679"),
    @(177, 46071.82291666666, 0, 0, "18.02.202680"),
    @(178, 46071.83333333334, 0, 0, "18.02.202681"),
    @(179, 46071.84375, 0, 0, "18.02.202682"),
    @(180, 46071.85416666666, 0, 0, "18.02.202683"),
    @(181, 46071.86458333334, 0, 0, "18.02.202684"),
    @(182, 46071.875, 0, 0, "18.02.202685"),
    @(183, 46071.88541666666, 0, 0, "18.02.202686"),
    @(184, 46071.89583333334, 0, 0, "18.02.202687"),
    @(185, 46071.90625, 0, 0, "18.02.202688"),
    @(186, 46071.91666666666, 0, 0, "18.02.202689"),
    @(187, 46071.92708333334, 0, 0, "18.02.202690"),
    @(188, 46071.9375, 0, 0, "18.02.202691"),
    @(189, 46071.94791666666, 0, 0, "18.02.202692"),
    @(190, 46071.95833333334, 0, 0, "18.02.202693"),
    @(191, 46071.96875, 0, 0, "18.02.202694"),
    @(192, 46071.97916666666, 0, 0, "18.02.202695"),
    @(193, 46071.98958333334, 0, 0, "18.02.202696")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
}

Write-Output "Updated $($data.Count) rows"
